# Applies trade #20 close update to live_trading_results workbook
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.27
$summary.Range("B4").Value = -0.73
$summary.Range("B5").Value = -0.73
$summary.Range("B6").Value = 20
$summary.Range("B8").Value = 11
$summary.Range("B9").Value = 25

# --- Strategy Status sheet ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.27
$status.Range("D4").Value = 20
$status.Range("E4").Value = -0.73
$status.Range("F4").Value = -0.73
$status.Range("G4").Value = 25

# --- New trade row data (Trade #20) ---
$newRow = @{
    A = 20
    B = "2026-02-17"
    C = "08:21:48"
    D = "MarketMaking"
    E = "DOWN"
    F = 0.41
    G = 0.27
    H = "CLOSED"
    I = -34.1463
    J = -0.14
    K = 99.27
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 0.13
}

function Add-TradeRow($ws, $rowNum, $data) {
    $ws.Range("A$rowNum").Value = $data.A

    # Force text format so the date/time-like strings aren't auto-converted to
    # date/time serial numbers, then restore the default "Normal" style so the
    # cell ends up with plain text and no explicit (non-default) style index.
    $ws.Range("B$rowNum").NumberFormat = "@"
    $ws.Range("B$rowNum").Value = $data.B
    $ws.Range("B$rowNum").Style = "Normal"

    $ws.Range("C$rowNum").NumberFormat = "@"
    $ws.Range("C$rowNum").Value = $data.C
    $ws.Range("C$rowNum").Style = "Normal"

    $ws.Range("D$rowNum").Value = $data.D
    $ws.Range("E$rowNum").Value = $data.E
    $ws.Range("F$rowNum").Value = $data.F
    $ws.Range("G$rowNum").Value = $data.G
    $ws.Range("H$rowNum").Value = $data.H
    $ws.Range("I$rowNum").Value = $data.I
    $ws.Range("J$rowNum").Value = $data.J
    $ws.Range("K$rowNum").Value = $data.K
    $ws.Range("L$rowNum").Value = $data.L
    $ws.Range("M$rowNum").Value = $data.M
    $ws.Range("N$rowNum").Value = $data.N
    $ws.Range("O$rowNum").Value = $data.O
    $ws.Range("P$rowNum").Value = $data.P
    $ws.Range("Q$rowNum").Value = $data.Q
}

# --- All Trades sheet ---
$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades 21 $newRow

# --- MarketMaking sheet ---
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking 21 $newRow
